$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 2.11624
$ws.Range("H2").Value = 6.34872
$ws.Range("I2").Value = 0.1897594766532197
$ws.Range("J2").Value = 0.1897594766532197
$ws.Range("M2").Value = 14.861848
$ws.Range("N2").Value = 44.585544
$ws.Range("O2").Value = 0.09055189482833943
$ws.Range("P2").Value = 0.09055189482833945
$ws.Range("Q2").Value = 31.45123721152
$ws.Range("R2").Value = 283.06113490368
$ws.Range("S2").Value = 0.01718308017258308
$ws.Range("T2").Value = 0.01718308017258309
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 2.11624
$ws.Range("H3").Value = 6.34872
$ws.Range("I3").Value = 0.1897594766532197
$ws.Range("J3").Value = 0.1897594766532197
$ws.Range("O3").Value = 0.1893562842131466
$ws.Range("P3").Value = 0.1893562842131466
$ws.Range("Q3").Value = 65.76879946653332
$ws.Range("R3").Value = 591.9191951988
$ws.Range("S3").Value = 0.03593214939328502
$ws.Range("T3").Value = 0.03593214939328503
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 2.11624
$ws.Range("H4").Value = 6.34872
$ws.Range("I4").Value = 0.1897594766532197
$ws.Range("J4").Value = 0.1897594766532197
$ws.Range("M4").Value = 18.10188466666667
$ws.Range("N4").Value = 54.305654
$ws.Range("O4").Value = 0.1102931450066459
$ws.Range("P4").Value = 0.1102931450066459
$ws.Range("Q4").Value = 38.30793240698667
$ws.Range("R4").Value = 344.77139166288
$ws.Range("S4").Value = 0.02092916947489879
$ws.Range("T4").Value = 0.02092916947489881
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 2.11624
$ws.Range("H5").Value = 6.34872
$ws.Range("I5").Value = 0.1897594766532197
$ws.Range("J5").Value = 0.1897594766532197
$ws.Range("M5").Value = 100.0833306666667
$ws.Range("N5").Value = 300.249992
$ws.Range("O5").Value = 0.609798675951868
$ws.Range("P5").Value = 0.6097986759518681
$ws.Range("Q5").Value = 211.8003476900266
$ws.Range("R5").Value = 1906.20312921024
$ws.Range("S5").Value = 0.1157150776124528
$ws.Range("T5").Value = 0.1157150776124528
$ws.Range("I6").Value = 0.6160274054778138
$ws.Range("J6").Value = 0.6160274054778138
$ws.Range("M6").Value = 14.861848
$ws.Range("N6").Value = 44.585544
$ws.Range("O6").Value = 0.09055189482833943
$ws.Range("P6").Value = 0.09055189482833945
$ws.Range("Q6").Value = 102.1020103986
$ws.Range("R6").Value = 918.9180935874
$ws.Range("S6").Value = 0.05578244883220181
$ws.Range("T6").Value = 0.05578244883220181
$ws.Range("I7").Value = 0.6160274054778138
$ws.Range("J7").Value = 0.6160274054778138
$ws.Range("O7").Value = 0.1893562842131466
$ws.Range("P7").Value = 0.1893562842131466
$ws.Range("S7").Value = 0.1166486604747442
$ws.Range("T7").Value = 0.1166486604747442
$ws.Range("I8").Value = 0.6160274054778138
$ws.Range("J8").Value = 0.6160274054778138
$ws.Range("M8").Value = 18.10188466666667
$ws.Range("N8").Value = 54.305654
$ws.Range("O8").Value = 0.1102931450066459
$ws.Range("P8").Value = 0.1102931450066459
$ws.Range("Q8").Value = 124.36130530135
$ws.Range("R8").Value = 1119.25174771215
$ws.Range("S8").Value = 0.06794359996043235
$ws.Range("T8").Value = 0.06794359996043238
$ws.Range("I9").Value = 0.6160274054778138
$ws.Range("J9").Value = 0.6160274054778138
$ws.Range("M9").Value = 100.0833306666667
$ws.Range("N9").Value = 300.249992
$ws.Range("O9").Value = 0.609798675951868
$ws.Range("P9").Value = 0.6097986759518681
$ws.Range("Q9").Value = 687.5799879297999
$ws.Range("R9").Value = 6188.2198913682
$ws.Range("S9").Value = 0.3756526962104353
$ws.Range("T9").Value = 0.3756526962104354
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 1.793503666666667
$ws.Range("H10").Value = 5.380511
$ws.Range("I10").Value = 0.1608202836929164
$ws.Range("J10").Value = 0.1608202836929164
$ws.Range("M10").Value = 14.861848
$ws.Range("N10").Value = 44.585544
$ws.Range("O10").Value = 0.09055189482833943
$ws.Range("P10").Value = 0.09055189482833945
$ws.Range("Q10").Value = 26.65477888144267
$ws.Range("R10").Value = 239.893009932984
$ws.Range("S10").Value = 0.01456258141522467
$ws.Range("T10").Value = 0.01456258141522468
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 1.793503666666667
$ws.Range("H11").Value = 5.380511
$ws.Range("I11").Value = 0.1608202836929164
$ws.Range("J11").Value = 0.1608202836929164
$ws.Range("O11").Value = 0.1893562842131466
$ws.Range("P11").Value = 0.1893562842131466
$ws.Range("Q11").Value = 55.73875505400722
$ws.Range("R11").Value = 501.648795486065
$ws.Range("S11").Value = 0.03045233134619473
$ws.Range("T11").Value = 0.03045233134619474
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 1.793503666666667
$ws.Range("H12").Value = 5.380511
$ws.Range("I12").Value = 0.1608202836929164
$ws.Range("J12").Value = 0.1608202836929164
$ws.Range("M12").Value = 18.10188466666667
$ws.Range("N12").Value = 54.305654
$ws.Range("O12").Value = 0.1102931450066459
$ws.Range("P12").Value = 0.1102931450066459
$ws.Range("Q12").Value = 32.46579652324378
$ws.Range("R12").Value = 292.192168709194
$ws.Range("S12").Value = 0.01773737486935275
$ws.Range("T12").Value = 0.01773737486935276
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 1.793503666666667
$ws.Range("H13").Value = 5.380511
$ws.Range("I13").Value = 0.1608202836929164
$ws.Range("J13").Value = 0.1608202836929164
$ws.Range("M13").Value = 100.0833306666667
$ws.Range("N13").Value = 300.249992
$ws.Range("O13").Value = 0.609798675951868
$ws.Range("P13").Value = 0.6097986759518681
$ws.Range("Q13").Value = 179.4998205228791
$ws.Range("R13").Value = 1615.498384705912
$ws.Range("S13").Value = 0.09806799606214418
$ws.Range("T13").Value = 0.09806799606214421
$ws.Range("E14").Value = 2.0
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.3724043333333333
$ws.Range("H14").Value = 1.117213
$ws.Range("I14").Value = 0.03339283417605023
$ws.Range("J14").Value = 0.03339283417605023
$ws.Range("M14").Value = 14.861848
$ws.Range("N14").Value = 44.585544
$ws.Range("O14").Value = 0.09055189482833943
$ws.Range("P14").Value = 0.09055189482833945
$ws.Range("Q14").Value = 5.534616596541333
$ws.Range("R14").Value = 49.811549368872
$ws.Range("S14").Value = 0.003023784408329879
$ws.Range("T14").Value = 0.003023784408329879
$ws.Range("E15").Value = 2.0
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.3724043333333333
$ws.Range("H15").Value = 1.117213
$ws.Range("I15").Value = 0.03339283417605023
$ws.Range("J15").Value = 0.03339283417605023
$ws.Range("O15").Value = 0.1893562842131466
$ws.Range("P15").Value = 0.1893562842131466
$ws.Range("Q15").Value = 11.57363338726611
$ws.Range("R15").Value = 104.162700485395
$ws.Range("S15").Value = 0.00632314299892264
$ws.Range("T15").Value = 0.006323142998922642
$ws.Range("E16").Value = 2.0
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.3724043333333333
$ws.Range("H16").Value = 1.117213
$ws.Range("I16").Value = 0.03339283417605023
$ws.Range("J16").Value = 0.03339283417605023
$ws.Range("M16").Value = 18.10188466666667
$ws.Range("N16").Value = 54.305654
$ws.Range("O16").Value = 0.1102931450066459
$ws.Range("P16").Value = 0.1102931450066459
$ws.Range("Q16").Value = 6.741220291366889
$ws.Range("R16").Value = 60.67098262230201
$ws.Range("S16").Value = 0.003683000701961988
$ws.Range("T16").Value = 0.003683000701961989
$ws.Range("E17").Value = 2.0
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.3724043333333333
$ws.Range("H17").Value = 1.117213
$ws.Range("I17").Value = 0.03339283417605023
$ws.Range("J17").Value = 0.03339283417605023
$ws.Range("M17").Value = 100.0833306666667
$ws.Range("N17").Value = 300.249992
$ws.Range("O17").Value = 0.609798675951868
$ws.Range("P17").Value = 0.6097986759518681
$ws.Range("Q17").Value = 37.27146603469955
$ws.Range("R17").Value = 335.443194312296
$ws.Range("S17").Value = 0.02036290606683571
$ws.Range("T17").Value = 0.02036290606683572
